$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.883.95"
$ws.Range("E2").Value = "  -4.12%  "
$ws.Range("D3").Value = "2.243.41"
$ws.Range("E3").Value = "  -5.43%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.34"
$ws.Range("E5").Value = "  -2.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.44"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -4.58%  "
$ws.Range("D9").Value = "2.247.58"
$ws.Range("E9").Value = "  -5.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0913"
$ws.Range("E10").Value = "  -7.03%  "
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.70"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.315"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "2.653.78"
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.14"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "53.859.75"
$ws.Range("E16").Value = "  -4.12%  "
$ws.Range("E17").Value = "  -3.05%  "
$ws.Range("D18").Value = "2.267.34"
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.97"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.58"
$ws.Range("E20").Value = "  -4.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "299.82"
$ws.Range("E21").Value = "  -2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.80"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.365"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.143"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.06"
$ws.Range("E28").Value = "  -3.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.43"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").Value = "0.0₃0687"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.59"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").Value = "  -0.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.06"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.40"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("E38").Value = "  +6.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.58"
$ws.Range("E39").Value = "  -5.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.68"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.367"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.37"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.30"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "122.26"
$ws.Range("E44").Value = "  -6.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.65"
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0877"
$ws.Range("E46").Value = "  -2.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.534"
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "236.34"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0470"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0202"
$ws.Range("E50").Value = "  -3.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.17"
$ws.Range("E51").Value = "  -4.08%  "
